$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) cells to text format so numeric-looking strings
# like "291.89" are preserved exactly as text, not coerced to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '22.446.03'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '1.571.97'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = '291.89'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '0.3724'
$ws.Range('E7').Value = '  -1.05%  '
$ws.Range('D8').Value = '49.95'
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('D9').Value = '0.3394'
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = '1.148'
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.07563'
$ws.Range('E11').Value = '  -1.21%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').Value = '6.029'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = '6.956'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').Value = '1.573.94'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').Value = '0.00001124'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '91.07'
$ws.Range('E18').Value = '  +1.08%  '
$ws.Range('D19').Value = '0.06760'
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '6.300'
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('D22').Value = '16.31'
$ws.Range('E22').Value = '  -2.90%  '
$ws.Range('D23').Value = '12.15'
$ws.Range('E23').Value = '  +1.04%  '
$ws.Range('D24').Value = '22.443.37'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '2.344'
$ws.Range('E25').Value = '  -2.24%  '
$ws.Range('D26').Value = '2.692'
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('D27').Value = '20.10'
$ws.Range('E27').Value = '  -0.49%  '
$ws.Range('D28').Value = '148.81'
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('D29').Value = '5.038'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '125.51'
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('D31').Value = '1.748.83'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('E32').Value = '  +6.96%  '
$ws.Range('D33').Value = '6.172'
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('E34').Value = '  -0.58%  '
$ws.Range('D35').Value = '9.853'
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('D36').Value = '0.08367'
$ws.Range('D37').Value = '0.02480'
$ws.Range('E37').Value = '  -2.16%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '0.2307'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '1.347'
$ws.Range('E39').Value = '  -3.39%  '
$ws.Range('D40').Value = '0.06522'
$ws.Range('D41').Value = '5.458'
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('D42').Value = '11.34'
$ws.Range('E42').Value = '  -1.20%  '
$ws.Range('D43').Value = '0.6228'
$ws.Range('E43').Value = '  -2.83%  '
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').Value = '14.07'
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('D46').Value = '3.804'
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.5818'
$ws.Range('E47').Value = '  -2.73%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '130.46'
$ws.Range('E48').Value = '  +3.85%  '
$ws.Range('D49').Value = '2.071'
$ws.Range('E49').Value = '  -0.99%  '
$ws.Range('E50').Value = '  -5.81%  '
$ws.Range('E51').Value = '  -0.26%  '

# Restore original (default/general) formatting on column D now that the
# text values have been written, so styling matches the source workbook.
$ws.Range("D2:D51").ClearFormats()
